$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 1.684223175048828
$ws.Range("B1").Value = 3.317033052444458
$ws.Range("C1").Value = 6.005094051361084
$ws.Range("D1").Value = 1.816447973251343
$ws.Range("E1").Value = 0.8975754976272583
